$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 15158181
$ws.Range("I51").Value = 30308528
$ws.Range("J51").Value = 7833.3335
$ws.Range("K51").Value = 30308528
$ws.Range("L51").Value = 7833.3335
$ws.Range("M51").Value = -30308044
$ws.Range("N51").Value = -8801.333500000001
$ws.Range("H74").Value = 3848.8572
$ws.Range("I74").Value = 3660.3333
$ws.Range("J74").Value = 4980
$ws.Range("K74").Value = 3660.3333
$ws.Range("L74").Value = 4980
$ws.Range("M74").Value = -2724.3333
$ws.Range("N74").Value = -6852
$ws.Range("H77").Value = 3848.8572
$ws.Range("I77").Value = 3660.3333
$ws.Range("J77").Value = 4980
$ws.Range("K77").Value = 18301.6665
$ws.Range("L77").Value = 24900
$ws.Range("M77").Value = -13621.6665
$ws.Range("N77").Value = -34260
$ws.Range("H98").Value = 1149.4242
$ws.Range("I98").Value = 981.28
$ws.Range("J98").Value = 1674.875
$ws.Range("K98").Value = 981.28
$ws.Range("L98").Value = 1674.875
$ws.Range("M98").Value = 516.72
$ws.Range("N98").Value = -4670.875
$ws.Range("H107").Value = 371.9091
$ws.Range("I107").Value = 349.6
$ws.Range("K107").Value = 349.6
$ws.Range("M107").Value = 1570.4
$ws.Range("H122").Value = 1149.4242
$ws.Range("I122").Value = 981.28
$ws.Range("J122").Value = 1674.875
$ws.Range("K122").Value = 2943.84
$ws.Range("L122").Value = 5024.625
$ws.Range("M122").Value = -493.8400000000001
$ws.Range("N122").Value = -9924.625
$ws.Range("H134").Value = 333364500
$ws.Range("J134").Value = 333364500
$ws.Range("L134").Value = 333364500
$ws.Range("N134").Value = -333374640
$ws.Range("H137").Value = 2136.361
$ws.Range("I137").Value = 2282.682
$ws.Range("J137").Value = 1906.4286
$ws.Range("K137").Value = 6848.045999999999
$ws.Range("L137").Value = 5719.2858
$ws.Range("M137").Value = -4298.045999999999
$ws.Range("N137").Value = -10819.2858

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7234.973
$ws.Range("I32").Value = 7622.719
$ws.Range("K32").Value = 7622.719
$ws.Range("M32").Value = -7335.719
$ws.Range("H45").Value = 2679.75
$ws.Range("I45").Value = 1681.7778
$ws.Range("K45").Value = 1681.7778
$ws.Range("M45").Value = -1304.7778
$ws.Range("H61").Value = 3714.7058
$ws.Range("I61").Value = 2200.5715
$ws.Range("J61").Value = 4774.6
$ws.Range("K61").Value = 2200.5715
$ws.Range("L61").Value = 4774.6
$ws.Range("M61").Value = -1988.5715
$ws.Range("N61").Value = -5198.6
$ws.Range("H86").Value = 19542.666
$ws.Range("J86").Value = 19542.666
$ws.Range("L86").Value = 19542.666
$ws.Range("N86").Value = -21914.666
$ws.Range("H89").Value = 19542.666
$ws.Range("J89").Value = 19542.666
$ws.Range("L89").Value = 58627.99800000001
$ws.Range("N89").Value = -70483.99800000001
$ws.Range("H92").Value = 35498.75
$ws.Range("J92").Value = 35498.75
$ws.Range("L92").Value = 35498.75
$ws.Range("N92").Value = -40490.75
$ws.Range("H122").Value = 1647.7391
$ws.Range("I122").Value = 1467.75
$ws.Range("J122").Value = 2059.1428
$ws.Range("K122").Value = 4403.25
$ws.Range("L122").Value = 6177.428400000001
$ws.Range("M122").Value = -1953.25
$ws.Range("N122").Value = -11077.4284
$ws.Range("H136").Value = 3714.7058
$ws.Range("I136").Value = 2200.5715
$ws.Range("J136").Value = 4774.6
$ws.Range("K136").Value = 6601.7145
$ws.Range("L136").Value = 14323.8
$ws.Range("M136").Value = -4051.7145
$ws.Range("N136").Value = -19423.8

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6538004
$ws.Range("I31").Value = 1553.0488
$ws.Range("J31").Value = 33337454
$ws.Range("K31").Value = 1553.0488
$ws.Range("L31").Value = 33337454
$ws.Range("M31").Value = -1258.0488
$ws.Range("N31").Value = -33338044
$ws.Range("H34").Value = 6538004
$ws.Range("I34").Value = 1553.0488
$ws.Range("J34").Value = 33337454
$ws.Range("K34").Value = 1553.0488
$ws.Range("L34").Value = 33337454
$ws.Range("M34").Value = -1351.0488
$ws.Range("N34").Value = -33337858
$ws.Range("H99").Value = 3081.4443
$ws.Range("I99").Value = 2096
$ws.Range("J99").Value = 5643.6
$ws.Range("K99").Value = 2096
$ws.Range("L99").Value = 5643.6
$ws.Range("M99").Value = -598
$ws.Range("N99").Value = -8639.6
$ws.Range("H126").Value = 3081.4443
$ws.Range("I126").Value = 2096
$ws.Range("J126").Value = 5643.6
$ws.Range("K126").Value = 6288
$ws.Range("L126").Value = 16930.8
$ws.Range("M126").Value = -3818
$ws.Range("N126").Value = -21870.8

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 3270
$ws.Range("J110").Value = 3937.5
$ws.Range("L110").Value = 11812.5
$ws.Range("N110").Value = -19992.5
$ws.Range("H120").Value = 10055.777
$ws.Range("I120").Value = 4821.6
$ws.Range("J120").Value = 12068.923
$ws.Range("K120").Value = 14464.8
$ws.Range("L120").Value = 36206.769
$ws.Range("M120").Value = -9626.800000000001
$ws.Range("N120").Value = -45882.769
$ws.Range("H124").Value = 19266.666
$ws.Range("I124").Value = 1800
$ws.Range("K124").Value = 5400
$ws.Range("M124").Value = -490
$ws.Range("H132").Value = 3504.9092
$ws.Range("I132").Value = 1850
$ws.Range("J132").Value = 4125.5
$ws.Range("K132").Value = 16650
$ws.Range("L132").Value = 37129.5
$ws.Range("M132").Value = -14120
$ws.Range("N132").Value = -42189.5

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2590.3
$ws.Range("I82").Value = 1200
$ws.Range("J82").Value = 3517.1667
$ws.Range("K82").Value = 1200
$ws.Range("L82").Value = 3517.1667
$ws.Range("M82").Value = -839
$ws.Range("N82").Value = -4239.1667
$ws.Range("H85").Value = 2590.3
$ws.Range("I85").Value = 1200
$ws.Range("J85").Value = 3517.1667
$ws.Range("K85").Value = 1200
$ws.Range("L85").Value = 3517.1667
$ws.Range("M85").Value = 48
$ws.Range("N85").Value = -6013.1667

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3662.5
$ws.Range("I81").Value = 4720.2
$ws.Range("J81").Value = 1899.6666
$ws.Range("K81").Value = 9440.4
$ws.Range("L81").Value = 3799.3332
$ws.Range("M81").Value = -8379.4
$ws.Range("N81").Value = -5921.3332
$ws.Range("H84").Value = 3662.5
$ws.Range("I84").Value = 4720.2
$ws.Range("J84").Value = 1899.6666
$ws.Range("K84").Value = 47202
$ws.Range("L84").Value = 18996.666
$ws.Range("M84").Value = -41898
$ws.Range("N84").Value = -29604.666
